# Updated symbol list on Sun Dec 18 12:44:10 UTC 2022 with GitHub Actions
#
# This mirrors the automated "refresh crypto prices" commit: a handful of
# Price cells (column D) got refreshed quotes, and two Volume(1h) label
# cells (column E) had their "Worstin24h" suffix move to a different coin.
#
# Column D cells store their numbers as literal text (t="inlineStr") rather
# than as real numeric cells, matching how the upstream openpyxl export
# writes them. A leading apostrophe forces Excel to keep the literal text
# (incl. trailing zeros) instead of re-parsing it as a float.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value  = "'5.526"
$ws.Range("D5").Value  = "'0.05628"
$ws.Range("D6").Value  = "'3.380"
$ws.Range("D7").Value  = "'6.471"
$ws.Range("D8").Value  = "'0.8051"
$ws.Range("D9").Value  = "'1.054"
$ws.Range("D10").Value = "'0.1423"
$ws.Range("D11").Value = "'0.07320"
$ws.Range("D12").Value = "'0.03192"
$ws.Range("D13").Value = "'0.02971"
$ws.Range("D14").Value = "'0.09264"
$ws.Range("D15").Value = "'0.001668"
$ws.Range("D16").Value = "'3.221"
$ws.Range("D17").Value = "'0.04703"
$ws.Range("D18").Value = "'0.0005878"
$ws.Range("D19").Value = "'0.006276"
$ws.Range("D20").Value = "'0.001060"
$ws.Range("D21").Value = "'0.004122"

$ws.Range("E23").Value = "22UpBotsUBXT"

$ws.Range("D24").Value = "'3.968"
$ws.Range("D25").Value = "'2.135"
$ws.Range("D26").Value = "'0.3272"

$ws.Range("D40").Value = "'0.04171"
$ws.Range("D41").Value = "'0.006878"

$ws.Range("D43").Value = "'0.1040"
$ws.Range("D44").Value = "'0.009890"
$ws.Range("D45").Value = "'0.00005648"

$ws.Range("D48").Value = "'0.02457"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

$ws.Range("D49").Value = "'0.00002102"
